$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.255.60'

$ws.Range('D3').Value = '1.871.07'
$ws.Range('E3').Value = '  +3.65%  '

$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.68'
$ws.Range('E5').Value = '  +1.51%  '

$ws.Range('E6').Value = '  +0.21%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5056'
$ws.Range('E7').Value = '  +1.23%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3923'
$ws.Range('E8').Value = '  +1.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09616'
$ws.Range('E9').Value = '  +1.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.145'
$ws.Range('E10').Value = '  +4.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.91'
$ws.Range('E11').Value = '  +1.14%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.478'
$ws.Range('E12').Value = '  +1.82%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.99'
$ws.Range('E13').Value = '  +2.58%  '

$ws.Range('D14').Value = '1.879.11'
$ws.Range('E14').Value = '  +3.67%  '

$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.002'
$ws.Range('E15').Value = '  +0.23%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.414'
$ws.Range('E16').Value = '  +2.61%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001129'
$ws.Range('E17').Value = '  +0.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.86'
$ws.Range('E18').Value = '  +0.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06600'
$ws.Range('E19').Value = '  +0.40%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.66'
$ws.Range('E20').Value = '  +3.16%  '

$ws.Range('E21').Value = '  +0.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.188'
$ws.Range('E22').Value = '  +4.21%  '

$ws.Range('D23').Value = '28.322.85'
$ws.Range('E23').Value = '  +1.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  +1.89%  '

$ws.Range('E25').Value = '  +2.84%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.567'
$ws.Range('E26').Value = '  +6.75%  '

$ws.Range('D27').Value = '2.093.66'
$ws.Range('E27').Value = '  +3.92%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.25'
$ws.Range('E28').Value = '  +3.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.60'
$ws.Range('E29').Value = '  +1.37%  '

$ws.Range('E30').Value = '  -0.11%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1065'
$ws.Range('E31').Value = '  -0.80%  '

$ws.Range('E32').Value = '  +1.17%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.626'
$ws.Range('E33').Value = '  +1.13%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.623'
$ws.Range('E34').Value = '  +0.43%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06752'
$ws.Range('E35').Value = '  -0.63%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.511'
$ws.Range('E36').Value = '  +6.79%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02404'
$ws.Range('E37').Value = '  +4.48%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2187'
$ws.Range('E38').Value = '  +2.25%  '

$ws.Range('E39').Value = '  +1.20%  '

$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.001'
$ws.Range('E40').Value = '  +1.53%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6346'
$ws.Range('E41').Value = '  +2.41%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.185'
$ws.Range('E42').Value = '  +3.49%  '

$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.55'
$ws.Range('E44').Value = '  +3.80%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5986'
$ws.Range('E45').Value = '  +2.14%  '

$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.272'
$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.659'
$ws.Range('E47').Value = '  -0.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.003'
$ws.Range('E48').Value = '  +2.78%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.11'
$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('E50').Value = '  +1.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06852'
$ws.Range('E51').Value = '  +1.62%  '
